# Add conditions for OBE1 rows & align stripplot by ID.
# Fills in the previously-"Unknown" Condition (col E) / Order Condition
# (col F) values for the OBE1 block (rows 2-19) with the correct
# condition codes, matching the pattern already used for later blocks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> @{ E = <condition code or $null>; F = <order condition code or $null> }
$updates = @{
    2  = @{ E = "IC" }
    3  = @{ E = "IC"; F = "I" }
    4  = @{ E = "IC"; F = "C" }
    5  = @{ E = "CI" }
    6  = @{ E = "IC"; F = "I" }
    7  = @{ E = "IC"; F = "C" }
    8  = @{ E = "CI" }
    9  = @{ E = "IC" }
    10 = @{ E = "CI" }
    11 = @{ E = "IC" }
    12 = @{ E = "CI" }
    13 = @{ E = "IC"; F = "I" }
    14 = @{ E = "IC"; F = "C" }
    15 = @{ E = "CI" }
    16 = @{ E = "IC" }
    18 = @{ E = "CI"; F = "C" }
    19 = @{ E = "CI"; F = "I" }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals.ContainsKey("E")) {
        $ws.Cells.Item($row, 5).Value = $vals["E"]
    }
    if ($vals.ContainsKey("F")) {
        $ws.Cells.Item($row, 6).Value = $vals["F"]
    }
}

# Active cell ended up on H9 when the editor finished working.
$ws.Range("H9").Select()
